# Updates the multiplication problems in the table to match the new
# randomly-generated set of equations (commit "Update master to output
# generated at 503736d").

$d = $word.ActiveDocument

# Mapping of old equation text -> new equation text. Every old value is
# unique within the document, so a MatchCase, non-wildcard, whole-document
# ReplaceAll is safe and will touch exactly one run each.
$replacements = @(
    @("773×6=", "753×8="),
    @("417×7=", "783×2="),
    @("688×7=", "141×6="),
    @("106×9=", "424×9="),
    @("571×9=", "727×6="),
    @("213×6=", "216×3="),
    @("648×2=", "152×3="),
    @("561×2=", "344×6="),
    @("246×3=", "357×4="),
    @("931×2=", "306×6="),
    @("615×2=", "477×2="),
    @("628×8=", "331×6="),
    @("237×9=", "891×7="),
    @("864×9=", "771×4="),
    @("584×7=", "631×9="),
    @("410×3=", "303×4="),
    @("128×3=", "861×2="),
    @("729×4=", "800×4="),
    @("521×2=", "341×2="),
    @("396×9=", "204×9="),
    @("518×6=", "786×7="),
    @("617×2=", "264×8="),
    @("438×5=", "762×9="),
    @("718×2=", "184×7="),
    @("797×7=", "339×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}
